# Update master to output generated at c8c62b6
$d = $word.ActiveDocument

# 1. Update the date/weekday heading at the top of the document.
$d.Content.Find.Execute("2025-08-27 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-08-28 Thursday", 2)

# 2. Update the division problems inside the practice table.
#    The table has 20 rows x 5 columns; only rows 1,5,9,13,17 (1-indexed)
#    hold the actual problems, the others are blank spacer rows.
$t = $d.Tables(1)

$rows = @(1, 5, 9, 13, 17)
$replacements = @(
    @("16÷5=", "22÷5=", "51÷6=", "59÷5=", "65÷4="),
    @("63÷8=", "69÷3=", "16÷3=", "98÷9=", "42÷9="),
    @("79÷3=", "92÷6=", "29÷2=", "81÷7=", "23÷2="),
    @("70÷8=", "26÷6=", "52÷2=", "39÷8=", "52÷5="),
    @("55÷2=", "89÷6=", "53÷7=", "47÷9=", "16÷5=")
)
$newValues = @(
    @("46÷7=", "35÷9=", "84÷8=", "59÷2=", "76÷7="),
    @("14÷2=", "45÷6=", "81÷2=", "75÷6=", "83÷5="),
    @("82÷5=", "31÷6=", "34÷3=", "49÷7=", "41÷8="),
    @("80÷9=", "97÷3=", "38÷2=", "80÷4=", "73÷2="),
    @("41÷7=", "89÷8=", "50÷2=", "45÷2=", "79÷8=")
)

for ($ri = 0; $ri -lt $rows.Length; $ri++) {
    $rowIndex = $rows[$ri]
    for ($ci = 1; $ci -le 5; $ci++) {
        $cell = $t.Cell($rowIndex, $ci)
        $r = $cell.Range
        # Trim the trailing end-of-cell mark so only the visible text is replaced
        # (this also keeps the replacement scoped to this single cell, since a
        # plain Find/Replace on the cell Range can otherwise match duplicate
        # problem text elsewhere in the table, e.g. "16÷5=" appears twice).
        $r.End = $r.End - 1
        $r.Text = $newValues[$ri][$ci - 1]
    }
}
